# "home page and fixtures"
# - rename the demo sheet from "ExcelModuleDemoToDoItem" to "DemoToDoItem"
# - update the active selection on that sheet to H42 (was A2)

$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ExcelModuleDemoToDoItem")
$ws.Name = "DemoToDoItem"

$ws.Activate() | Out-Null
$ws.Range("H42").Select() | Out-Null
